# Femacal de La Calera - Haba: insert a new weekly record.
#
# The source data table (row 1 = headers, data starts row 2) is ordered by
# date descending in a way that this edit inserts one new record right
# after the existing row 62 (Sheet1!A62:R62), pushing every following row
# down by one. The new row that lands at position 63 carries the values
# that used to live in row 62, and row 62 itself is updated with a new
# sampling date / volume (everything else on that record stays the same).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the current row 62 before it gets touched - this is the data
# that needs to be duplicated down into the freshly inserted row 63.
# NOTE: reads must call the getter as a method, `.Value()`, not the bare
# property `.Value` - the latter yields the member descriptor itself here.
$carry = @(
    $ws.Cells.Item(62, 1).Value(),
    $ws.Cells.Item(62, 2).Value(),
    $ws.Cells.Item(62, 3).Value(),
    $ws.Cells.Item(62, 4).Value(),
    $ws.Cells.Item(62, 5).Value(),
    $ws.Cells.Item(62, 6).Value(),
    $ws.Cells.Item(62, 7).Value(),
    $ws.Cells.Item(62, 8).Value(),
    $ws.Cells.Item(62, 9).Value(),
    $ws.Cells.Item(62, 10).Value(),
    $ws.Cells.Item(62, 11).Value(),
    $ws.Cells.Item(62, 12).Value(),
    $ws.Cells.Item(62, 13).Value(),
    $ws.Cells.Item(62, 14).Value(),
    $ws.Cells.Item(62, 15).Value(),
    $ws.Cells.Item(62, 16).Value(),
    $ws.Cells.Item(62, 17).Value(),
    $ws.Cells.Item(62, 18).Value()
)

# Push rows 63:170 down to 64:171, opening up a blank row at 63.
$ws.Rows.Item(63).Insert()

# Fill the newly opened row 63 with the values that used to be in row 62.
$ws.Cells.Item(63, 1).Value = $carry[0]
$ws.Cells.Item(63, 2).Value = $carry[1]
$ws.Cells.Item(63, 3).Value = $carry[2]
$ws.Cells.Item(63, 4).Value = $carry[3]
$ws.Cells.Item(63, 5).Value = $carry[4]
$ws.Cells.Item(63, 6).Value = $carry[5]
$ws.Cells.Item(63, 7).Value = $carry[6]
$ws.Cells.Item(63, 8).Value = $carry[7]
$ws.Cells.Item(63, 9).Value = $carry[8]
$ws.Cells.Item(63, 10).Value = $carry[9]
$ws.Cells.Item(63, 11).Value = $carry[10]
$ws.Cells.Item(63, 12).Value = $carry[11]
$ws.Cells.Item(63, 13).Value = $carry[12]
$ws.Cells.Item(63, 14).Value = $carry[13]
$ws.Cells.Item(63, 15).Value = $carry[14]
$ws.Cells.Item(63, 16).Value = $carry[15]
$ws.Cells.Item(63, 17).Value = $carry[16]
$ws.Cells.Item(63, 18).Value = $carry[17]

# Update row 62 in place with the new record's date (D) and volume (J);
# every other field on that row is unchanged.
$ws.Cells.Item(62, 4).Value = 44775
$ws.Cells.Item(62, 10).Value = 45
